# Applies the "II Implementations" sheet updates: exported diff. II
# implementation data (sizes/times per encoding scheme) for the two
# result tables, plus a few label/annotation tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("II Implementations")

# ---------------------------------------------------------------------
# Small text / label edits
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Time Results (sec) * 10^6 = nano sec"
$ws.Range("K3").Value = "Refer to original file if something is not clear here"
$ws.Range("H6").Value = "WT similarity"
$ws.Range("A29").Value = "Memory (MB)"

# ---------------------------------------------------------------------
# First results table: header row 6 (A6:O6) already present; data rows
# 7-16, columns B:O (column A already holds the dataset-length values).
# ---------------------------------------------------------------------
$table1 = @(
    @(2.6,17500,504,19.4,10.1,20.2,261,87.6264,2.78,5.07,44.8962,5.21,10.4,20.2),
    @(5.1,29600,1010,78.9,22.8,45.4,484,148.563,5.06,32,79.0068,10.6,21.8,42.2),
    @(7.7,75300,2040,73.5,34.7,84.3,1170,397.913,7.77,7.56,71.2664,16.7,33.7,65.9),
    @(11,63000,2300,186,58.3,105,1120,348.627,7.05,62.2,253.421,22.2,45,86.5),
    @(14,103000,3010,139,58.5,128,1790,548.898,18.9,48.1,294.651,26.6,53.5,105),
    @(17,138000,4200,99.4,75.5,164,2350,884.906,18.7,18.7,389.859,34.7,68.8,135),
    @(20,148000,4210,345,88.1,187,2440,885.567,13.1,111,383.446,38.1,77.5,152),
    @(23,95300,4090,522,95.3,190,1750,587.96,66.5,178,229.129,44,88.1,170),
    @(25,118000,5100,361,110,225,2270,656.274,54.3,143,456.243,49.3,99.4,195),
    @(28,199000,6280,286,122,277,3280,1164.77,65.2,105,439.857,55.7,129,231)
)

# Second results table: header row 31 already present; data rows 32-41,
# columns B:O (ratios relative to the first table).
$table2 = @(
    @(4.89,0.928411,0.93,2.78,3.06,2.26,1.6,0.509718,3.6,2.76,1.69208,2.31,1.62,1.18),
    @(9.77,1.916049,1.8,5.53,6.12,4.52,3.36,0.998991,7.25,5.49,3.36832,4.61,3.23,2.34),
    @(14.65,2.908364,2.69,8.29,9.18,6.78,5.32,1.48751,10.89,8.21,5.01489,6.91,4.83,3.51),
    @(19.54,3.854686,3.57,11.05,12.23,9.03,7.09,1.97666,14.51,11,6.67939,9.21,6.44,4.68),
    @(24.43,4.807434,4.45,13.81,15.28,11.29,8.84,2.4656,18.11,13.77,8.3454,11.5,8.05,5.84),
    @(29.35,5.73892,5.35,16.6,18.33,13.54,11.19,2.95438,21.73,16.49,10.0678,13.8,9.66,7.01),
    @(34.28,6.651574,6.24,19.37,21.4,15.8,13.07,3.44338,25.28,19.24,11.6451,16.1,11.26,8.18),
    @(39.26,7.590143,7.14,22.18,24.45,18.06,14.93,3.93453,28.97,21.99,13.3462,18.41,12.88,9.36),
    @(44.16,8.480882,8.03,24.95,27.5,20.31,16.77,4.42175,32.59,24.74,14.8824,20.71,14.49,10.52),
    @(49.07,9.396931,8.92,27.72,30.54,22.56,18.64,4.91021,36.16,27.48,16.6276,23,16.09,11.69)
)

# Column C of table 1 gets a slightly different number format (xf with an
# (empty) alignment block) than the rest of the numeric cells; reproduce
# that by formatting column C separately from columns B,D:O.
for ($i = 0; $i -lt $table1.Length; $i++) {
    $r = 7 + $i
    $row = $table1[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = 2 + $j   # column B = 2
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}

for ($i = 0; $i -lt $table2.Length; $i++) {
    $r = 32 + $i
    $row = $table2[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = 2 + $j   # column B = 2
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}

# ---------------------------------------------------------------------
# Number formats
# ---------------------------------------------------------------------
# Table 1: integer-looking format "0" across B7:O16, with column C on its
# own xf (numFmtId 1 as well, but with an alignment attribute recorded).
$ws.Range("B7:O16").NumberFormat = "0"
$ws.Range("C7:C16").NumberFormat = "0"
$ws.Range("C7:C16").HorizontalAlignment = -4108

# Table 2: one-decimal format "0.0" across B32:O41.
$ws.Range("B32:O41").NumberFormat = "0.0"

# ---------------------------------------------------------------------
# Column widths picked up by Excel's "best fit" for the newly-populated
# columns C and H.
# ---------------------------------------------------------------------
$ws.Range("C1").ColumnWidth = 11.1640625
$ws.Range("H1").ColumnWidth = 12.1640625

# ---------------------------------------------------------------------
# Selection left where the author's cursor ended up after entering data.
# ---------------------------------------------------------------------
$ws.Range("L37").Select()
